$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '29.471.50'
Set-TextValue $ws 'E2' '  -0.15%  '
Set-TextValue $ws 'D3' '1.904.23'
Set-TextValue $ws 'E3' '  -0.57%  '
Set-TextValue $ws 'E4' '  +0.11%  '
Set-TextValue $ws 'D5' '325.66'
Set-TextValue $ws 'E5' '  -2.52%  '
Set-TextValue $ws 'D6' '1.003'
Set-TextValue $ws 'E6' '  +0.13%  '
Set-TextValue $ws 'D7' '0.4800'
Set-TextValue $ws 'E7' '  +2.59%  '
Set-TextValue $ws 'D8' '0.4071'
Set-TextValue $ws 'E8' '  -0.66%  '
Set-TextValue $ws 'D9' '0.08073'
Set-TextValue $ws 'E9' '  +0.42%  '
Set-TextValue $ws 'E10' '  -1.04%  '
Set-TextValue $ws 'D11' '23.38'
Set-TextValue $ws 'E11' '  +4.59%  '
Set-TextValue $ws 'D12' '1.913.50'
Set-TextValue $ws 'E12' '  +1.12%  '
Set-TextValue $ws 'D13' '5.956'
Set-TextValue $ws 'E13' '  -0.54%  '
Set-TextValue $ws 'D14' '7.081'
Set-TextValue $ws 'E14' '  -1.34%  '
Set-TextValue $ws 'D15' '90.10'
Set-TextValue $ws 'E15' '  +0.22%  '
Set-TextValue $ws 'E16' '  +0.16%  '
Set-TextValue $ws 'E17' '  +1.51%  '
Set-TextValue $ws 'E18' '  -0.14%  '
Set-TextValue $ws 'D19' '17.64'
Set-TextValue $ws 'E19' '  -0.99%  '
Set-TextValue $ws 'D20' '1.002'
Set-TextValue $ws 'E20' '  +0.19%  '
Set-TextValue $ws 'D21' '29.488.78'
Set-TextValue $ws 'E21' '  -0.04%  '
Set-TextValue $ws 'D22' '5.543'
Set-TextValue $ws 'E22' '  -0.67%  '
Set-TextValue $ws 'D23' '11.78'
Set-TextValue $ws 'E23' '  +1.97%  '
Set-TextValue $ws 'E24' '  -1.96%  '
Set-TextValue $ws 'D25' '2.142.68'
Set-TextValue $ws 'E25' '  -0.35%  '
Set-TextValue $ws 'D26' '154.57'
Set-TextValue $ws 'E26' '  -0.62%  '
Set-TextValue $ws 'D27' '19.84'
Set-TextValue $ws 'E27' '  -0.11%  '
Set-TextValue $ws 'D28' '6.107'
Set-TextValue $ws 'E28' '  +6.00%  '
Set-TextValue $ws 'D29' '2.097'
Set-TextValue $ws 'E29' '  -2.06%  '
Set-TextValue $ws 'D30' '118.34'
Set-TextValue $ws 'D31' '1.034'
Set-TextValue $ws 'E31' '  -3.22%  '
Set-TextValue $ws 'D32' '0.09512'
Set-TextValue $ws 'E32' '  +0.55%  '
Set-TextValue $ws 'D33' '5.521'
Set-TextValue $ws 'E34' '  -2.27%  '
Set-TextValue $ws 'E35' '  -0.99%  '
Set-TextValue $ws 'E36' '  -0.56%  '
Set-TextValue $ws 'D37' '0.06078'
Set-TextValue $ws 'E37' '  -0.72%  '
Set-TextValue $ws 'D38' '1.177'
Set-TextValue $ws 'E38' '  -0.38%  '
Set-TextValue $ws 'D39' '0.5903'
Set-TextValue $ws 'E39' '  +0.06%  '
Set-TextValue $ws 'D40' '7.939'
Set-TextValue $ws 'E40' '  -5.60%  '
Set-TextValue $ws 'D41' '0.1845'
Set-TextValue $ws 'E41' '  +0.16%  '
Set-TextValue $ws 'D42' '10.22'
Set-TextValue $ws 'E42' '  -0.10%  '
Set-TextValue $ws 'D43' '1.285'
Set-TextValue $ws 'E43' '  +1.27%  '
Set-TextValue $ws 'B44' 'RenderToken'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D44' '2.401'
Set-TextValue $ws 'E44' '  +1.19%  '
Set-TextValue $ws 'B45' 'Cronos'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D45' '0.07795'
Set-TextValue $ws 'E45' '  +3.85%  '
Set-TextValue $ws 'D46' '12.31'
Set-TextValue $ws 'E46' '  +0.73%  '
Set-TextValue $ws 'D47' '0.5537'
Set-TextValue $ws 'E47' '  -0.69%  '
Set-TextValue $ws 'D48' '1.925'
Set-TextValue $ws 'E48' '  -0.26%  '
Set-TextValue $ws 'D49' '114.10'
Set-TextValue $ws 'E49' '  +0.50%  '
Set-TextValue $ws 'B50' 'Aave'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D50' '72.43'
Set-TextValue $ws 'E50' '  +1.09%  '
Set-TextValue $ws 'B51' 'WOONetwork'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue $ws 'D51' '0.2941'
Set-TextValue $ws 'E51' '  -1.54%  '

Write-Output "Applied 97 cell updates"
